$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from the row above so the new date cell matches the
# existing date column style (numFmtId 14) instead of creating a new one.
$ws.Range("A51").Copy()
$ws.Range("A52").PasteSpecial(-4122)  # xlPasteFormats

# New row of parkrun data: date=2023-05-20 (serial 45066), parkrun_no=51, parkrun="colney", time=28.09
$ws.Cells.Item(52, 1).Value = "5/20/2023"
$ws.Cells.Item(52, 2).Value = 51
$ws.Cells.Item(52, 3).Value = "colney"
$ws.Cells.Item(52, 4).Value = 28.09

$ws.Range("E52").Select()
